$wb = $excel.ActiveWorkbook
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newWs.Name = "TestFreeze"
$newWs.Range("D33").Select()
$excel.ActiveWindow.SplitColumn = 3
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true
Write-Output "frozen"
